# edit.ps1 — apply the "affix list" revision:
#   * wrap previously-unmarked single-word runs in w:proofErr spellStart/spellEnd
#     (splitting the trailing-space run out into its own run where needed)
#   * drop the stray _GoBack bookmark that sat in the "Not " cell
#   * collapse the trailing "Things to decide: / Ize / Ty " paragraphs into a
#     single paragraph holding just the _GoBack bookmark + a space

$d = $word.ActiveDocument

function Wrap-BodyXml($inner) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        $inner + '<w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$transforms = @(
    ,@(10, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ial</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(11, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ic</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(15, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ous</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(16, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Eous</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(20, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ise</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>')
    ,@(21, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Acy</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(22, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>A</w:t></w:r><w:r><w:t>ry</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(24, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ious</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(25, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ity</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(26, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>En</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(27, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ful</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(28, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ive</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(29, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ative</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(30, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Itive</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(32, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ence</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ance</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(33, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ish</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>')
    ,@(35, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>I</w:t></w:r><w:r><w:t>ty</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(38, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>dom</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(43, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (present participle)</w:t></w:r></w:p>')
    ,@(47, '<w:p><w:r><w:t xml:space="preserve">Not </w:t></w:r></w:p>')
    ,@(55, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Im</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>')
    ,@(67, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Er</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(75, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Er</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(79, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ment</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(81, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>dess</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(82, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ess</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(83, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ist</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>')
    ,@(84, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ie</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(87, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>ster</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(91, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Tion</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(92, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:lastRenderedPageBreak/><w:t>Ment</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(93, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ation</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(94, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:lastRenderedPageBreak/><w:t>Ition</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(96, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:lastRenderedPageBreak/><w:t>Ible</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
    ,@(97, '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ize</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
)


# Apply every single-paragraph rewrite. Indices are 1-based positions into
# $d.Paragraphs as they exist in the (still unmodified-in-count) document —
# none of these per-paragraph rewrites add or remove paragraphs, so the
# numbering stays stable while we walk the list.
$paragraphs = $d.Paragraphs
foreach ($t in $transforms) {
    $idx = $t[0]
    $xml = $t[1]
    $p = $paragraphs.Item($idx)
    $p.Range.InsertXML((Wrap-BodyXml $xml))
}

# Collapse the final three paragraphs ("Things to decide:", "Ize", "Ty ")
# into one paragraph carrying the _GoBack bookmark + a single space.
$paragraphs2 = $d.Paragraphs
$n = $paragraphs2.Count
$pFirst = $paragraphs2.Item($n - 2)
$pLast = $paragraphs2.Item($n)
$tailRange = $d.Range($pFirst.Range.Start, $pLast.Range.End)
$tailInner = '<w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$tailRange.InsertXML((Wrap-BodyXml $tailInner))

Write-Output "done"
